$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "42.694.58"
$ws.Range("E2").Value = "  -0.04%  "
# Row 3
$ws.Range("D3").Value = "2.537.51"
$ws.Range("E3").Value = "  -0.39%  "
# Row 4
$ws.Range("E4").Value = "  +0.00%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "312.51"
$ws.Range("E5").Value = "  +0.62%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "100.17"
$ws.Range("E6").Value = "  +0.22%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.566"
$ws.Range("E7").Value = "  -0.87%  "
# Row 8
$ws.Range("E8").Value = "  +0.02%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.522"
$ws.Range("E9").Value = "  -2.10%  "
# Row 10
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "35.55"
$ws.Range("E10").Value = "  -1.17%  "
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0806"
$ws.Range("E11").Value = "  -0.23%  "
# Row 12
$ws.Range("B12").Value = "Polkadot"
$ws.Range("C12").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "7.29"
$ws.Range("E12").Value = "  -1.55%  "
# Row 13
$ws.Range("B13").Value = "TRON"
$ws.Range("C13").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.110"
$ws.Range("E13").Value = "  +0.58%  "
# Row 14
$ws.Range("D14").Value = "2.929.73"
$ws.Range("E14").Value = "  -0.19%  "
# Row 15
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").Value = "2.564.17"
$ws.Range("E15").Value = "  +0.07%  "
# Row 16
$ws.Range("B16").Value = "Chainlink"
$ws.Range("C16").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "15.43"
$ws.Range("E16").Value = "  -4.08%  "
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.814"
$ws.Range("E17").Value = "  -3.27%  "
# Row 18
$ws.Range("D18").Value = "42.701.65"
$ws.Range("E18").Value = "  -0.03%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.68"
$ws.Range("E19").Value = "  -1.30%  "
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "12.34"
$ws.Range("E20").Value = "  -0.33%  "
# Row 21
$ws.Range("E21").Value = "  -0.87%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "69.73"
$ws.Range("E22").Value = "  +0.31%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "242.85"
$ws.Range("E23").Value = "  -2.35%  "
# Row 24
$ws.Range("E24").Value = "  -1.91%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.02"
$ws.Range("E25").Value = "  -2.07%  "
# Row 27
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.61"
$ws.Range("E27").Value = "  -3.77%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.30"
$ws.Range("E28").Value = "  -3.15%  "
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "10.12"
$ws.Range("E29").Value = "  -0.12%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "38.30"
$ws.Range("E30").Value = "  -5.03%  "
# Row 31
$ws.Range("E31").Value = "  +2.49%  "
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "157.36"
$ws.Range("E32").Value = "  -0.22%  "
# Row 33
$ws.Range("E33").Value = "  +1.69%  "
# Row 34
$ws.Range("E34").Value = "  -1.92%  "
# Row 35
$ws.Range("B35").Value = "ApeXProtocol"
$ws.Range("C35").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.59"
$ws.Range("E35").Value = "  -0.75%  "
# Row 36
$ws.Range("B36").Value = "LidoDAOToken"
$ws.Range("C36").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.17"
$ws.Range("E36").Value = "  -4.11%  "
# Row 37
$ws.Range("B37").Value = "ARBITRUM"
$ws.Range("C37").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.96"
$ws.Range("E37").Value = "  -5.83%  "
# Row 38
$ws.Range("B38").Value = "Celestia"
$ws.Range("C38").Value = "https://coinranking.com/coin/YQcD0lBl7+celestia-tia"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "17.72"
$ws.Range("E38").Value = "  -2.83%  "
# Row 39
$ws.Range("E39").Value = "  -0.97%  "
# Row 40
$ws.Range("E40").Value = "  -0.33%  "
# Row 41
$ws.Range("E41").Value = "  -1.62%  "
# Row 42
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "21.78"
$ws.Range("E42").Value = "  -3.47%  "
# Row 43
$ws.Range("B43").Value = "NEARProtocol"
$ws.Range("C43").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "3.36"
$ws.Range("E43").Value = "  +2.55%  "
# Row 44
$ws.Range("B44").Value = "FirstDigitalUSD"
$ws.Range("C44").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.00"
$ws.Range("E44").Value = "  +0.22%  "
# Row 45
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.0299"
$ws.Range("E45").Value = "  -0.66%  "
# Row 46
$ws.Range("D46").Value = "1.992.01"
$ws.Range("E46").Value = "  +0.22%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "9.16"
$ws.Range("E47").Value = "  +1.41%  "
# Row 48
$ws.Range("D48").Value = "2.780.57"
$ws.Range("E48").Value = "  -0.04%  "
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.191"
$ws.Range("E49").Value = "  -1.45%  "
# Row 50
$ws.Range("E50").Value = "  -2.25%  "
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "72.29"
$ws.Range("E51").Value = "  -1.53%  "
